$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the header cell values first
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Then copy the formatting only from the existing header (H1) so the
# values we just wrote are preserved.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

# Add new data cells I2 and J2
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8
